$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bibi Cell Mundi
$ws.Range("O2").Value = 3636.9
$ws.Range("AG2").Value = 167649.86

# Row 3 - Bibi Cell Vieiralves
$ws.Range("O3").Value = 14071
$ws.Range("AG3").Value = 75525.8

# Row 4 - Bibi Cell Ponta Negra
$ws.Range("O4").Value = 4390.5
$ws.Range("P4").Value = 1481.42
$ws.Range("AG4").Value = 47569.69

# Row 5 - Bibi Cell Manauara
$ws.Range("N5").Value = 2493
$ws.Range("O5").Value = 5411
$ws.Range("P5").Value = 3140
$ws.Range("AG5").Value = 46405.2

# Row 6 - total
$ws.Range("N6").Value = 22435.29
$ws.Range("O6").Value = 27509.4
$ws.Range("P6").Value = 4621.42
$ws.Range("AG6").Value = 337150.55
